$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1839.8625
$ws.Cells.Item(40, 9).Value = 1874.091
$ws.Cells.Item(40, 10).Value = 1678.5
$ws.Cells.Item(40, 11).Value = 1874.091
$ws.Cells.Item(40, 12).Value = 1678.5
$ws.Cells.Item(40, 13).Value = -1699.091
$ws.Cells.Item(40, 14).Value = -2028.5

$ws.Cells.Item(62, 8).Value = 9525102
$ws.Cells.Item(62, 9).Value = 12346852
$ws.Cells.Item(62, 10).Value = 1697
$ws.Cells.Item(62, 11).Value = 12346852
$ws.Cells.Item(62, 12).Value = 1697
$ws.Cells.Item(62, 13).Value = -12346228
$ws.Cells.Item(62, 14).Value = -2945

$ws.Cells.Item(65, 8).Value = 9525102
$ws.Cells.Item(65, 9).Value = 12346852
$ws.Cells.Item(65, 10).Value = 1697
$ws.Cells.Item(65, 11).Value = 61734260
$ws.Cells.Item(65, 12).Value = 8485
$ws.Cells.Item(65, 13).Value = -61731140
$ws.Cells.Item(65, 14).Value = -14725

$ws.Cells.Item(100, 8).Value = 16668958
$ws.Cells.Item(100, 9).Value = 27779496
$ws.Cells.Item(100, 10).Value = 3150
$ws.Cells.Item(100, 11).Value = 27779496
$ws.Cells.Item(100, 12).Value = 3150
$ws.Cells.Item(100, 13).Value = -27778955
$ws.Cells.Item(100, 14).Value = -4232

$ws.Cells.Item(116, 8).Value = 6929.591
$ws.Cells.Item(116, 9).Value = 10004.583
$ws.Cells.Item(116, 10).Value = 3239.6
$ws.Cells.Item(116, 11).Value = 10004.583
$ws.Cells.Item(116, 12).Value = 3239.6
$ws.Cells.Item(116, 13).Value = -6562.583000000001
$ws.Cells.Item(116, 14).Value = -10123.6

$ws.Cells.Item(132, 8).Value = 2402.8948
$ws.Cells.Item(132, 9).Value = 2089.2144
$ws.Cells.Item(132, 10).Value = 3281.2
$ws.Cells.Item(132, 11).Value = 6267.6432
$ws.Cells.Item(132, 12).Value = 9843.599999999999
$ws.Cells.Item(132, 13).Value = -3737.6432
$ws.Cells.Item(132, 14).Value = -14903.6

$ws.Cells.Item(137, 8).Value = 1831.2821
$ws.Cells.Item(137, 9).Value = 1494.8846
$ws.Cells.Item(137, 10).Value = 2504.077
$ws.Cells.Item(137, 11).Value = 4484.6538
$ws.Cells.Item(137, 12).Value = 7512.231000000001
$ws.Cells.Item(137, 13).Value = -1934.6538
$ws.Cells.Item(137, 14).Value = -12612.231

$ws.Cells.Item(138, 8).Value = 1568.65
$ws.Cells.Item(138, 9).Value = 811.2
$ws.Cells.Item(138, 10).Value = 2188.3818
$ws.Cells.Item(138, 11).Value = 2433.6
$ws.Cells.Item(138, 12).Value = 6565.1454
$ws.Cells.Item(138, 13).Value = 2706.4
$ws.Cells.Item(138, 14).Value = -16845.1454

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4588
$ws.Cells.Item(32, 9).Value = 3438.8076
$ws.Cells.Item(32, 10).Value = 8662.409
$ws.Cells.Item(32, 11).Value = 3438.8076
$ws.Cells.Item(32, 12).Value = 8662.409
$ws.Cells.Item(32, 13).Value = -3151.8076
$ws.Cells.Item(32, 14).Value = -9236.409

$ws.Cells.Item(45, 8).Value = 4816.375
$ws.Cells.Item(45, 9).Value = 5004.0435
$ws.Cells.Item(45, 10).Value = 500
$ws.Cells.Item(45, 11).Value = 5004.0435
$ws.Cells.Item(45, 12).Value = 500
$ws.Cells.Item(45, 13).Value = -4627.0435
$ws.Cells.Item(45, 14).Value = -1254

$ws.Cells.Item(74, 8).Value = 1206.6595
$ws.Cells.Item(74, 9).Value = 1048.7949
$ws.Cells.Item(74, 10).Value = 1976.25
$ws.Cells.Item(74, 11).Value = 1048.7949
$ws.Cells.Item(74, 12).Value = 1976.25
$ws.Cells.Item(74, 13).Value = -174.7949000000001
$ws.Cells.Item(74, 14).Value = -3724.25

$ws.Cells.Item(77, 8).Value = 1206.6595
$ws.Cells.Item(77, 9).Value = 1048.7949
$ws.Cells.Item(77, 10).Value = 1976.25
$ws.Cells.Item(77, 11).Value = 5243.9745
$ws.Cells.Item(77, 12).Value = 9881.25
$ws.Cells.Item(77, 13).Value = -875.9745000000003
$ws.Cells.Item(77, 14).Value = -18617.25

$ws.Cells.Item(122, 8).Value = 1603342
$ws.Cells.Item(122, 9).Value = 1973096.9
$ws.Cells.Item(122, 11).Value = 5919290.699999999
$ws.Cells.Item(122, 13).Value = -5916840.699999999

$ws.Cells.Item(132, 8).Value = 2126.9265
$ws.Cells.Item(132, 9).Value = 1400.6
$ws.Cells.Item(132, 10).Value = 5199.846
$ws.Cells.Item(132, 11).Value = 4201.799999999999
$ws.Cells.Item(132, 12).Value = 15599.538
$ws.Cells.Item(132, 13).Value = -1671.799999999999
$ws.Cells.Item(132, 14).Value = -20659.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 7413.1113
$ws.Cells.Item(105, 9).Value = 11380.65
$ws.Cells.Item(105, 10).Value = 2453.6875
$ws.Cells.Item(105, 11).Value = 11380.65
$ws.Cells.Item(105, 12).Value = 2453.6875
$ws.Cells.Item(105, 13).Value = -9633.65
$ws.Cells.Item(105, 14).Value = -5947.6875

$ws.Cells.Item(134, 8).Value = 3060.817
$ws.Cells.Item(134, 9).Value = 3936.8108
$ws.Cells.Item(134, 10).Value = 2107.5293
$ws.Cells.Item(134, 11).Value = 11810.4324
$ws.Cells.Item(134, 12).Value = 6322.5879
$ws.Cells.Item(134, 13).Value = -9275.432400000002
$ws.Cells.Item(134, 14).Value = -11392.5879

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 209818.23
$ws.Cells.Item(31, 9).Value = 1498.3137
$ws.Cells.Item(31, 10).Value = 589258.0600000001
$ws.Cells.Item(31, 11).Value = 1498.3137
$ws.Cells.Item(31, 12).Value = 589258.0600000001
$ws.Cells.Item(31, 13).Value = -1203.3137
$ws.Cells.Item(31, 14).Value = -589848.0600000001

$ws.Cells.Item(34, 8).Value = 209818.23
$ws.Cells.Item(34, 9).Value = 1498.3137
$ws.Cells.Item(34, 10).Value = 589258.0600000001
$ws.Cells.Item(34, 11).Value = 1498.3137
$ws.Cells.Item(34, 12).Value = 589258.0600000001
$ws.Cells.Item(34, 13).Value = -1296.3137
$ws.Cells.Item(34, 14).Value = -589662.0600000001

$ws.Cells.Item(99, 8).Value = 5961249
$ws.Cells.Item(99, 9).Value = 8889
$ws.Cells.Item(99, 11).Value = 8889
$ws.Cells.Item(99, 13).Value = -7391

$ws.Cells.Item(126, 8).Value = 5961249
$ws.Cells.Item(126, 9).Value = 8889
$ws.Cells.Item(126, 11).Value = 26667
$ws.Cells.Item(126, 13).Value = -24197

$ws.Cells.Item(132, 8).Value = 1489.939
$ws.Cells.Item(132, 9).Value = 1117.5
$ws.Cells.Item(132, 10).Value = 3026.25
$ws.Cells.Item(132, 11).Value = 3352.5
$ws.Cells.Item(132, 12).Value = 9078.75
$ws.Cells.Item(132, 13).Value = -822.5
$ws.Cells.Item(132, 14).Value = -14138.75

$ws.Cells.Item(134, 8).Value = 1938.5518
$ws.Cells.Item(134, 9).Value = 2244.4866
$ws.Cells.Item(134, 10).Value = 1399.5238
$ws.Cells.Item(134, 11).Value = 6733.459800000001
$ws.Cells.Item(134, 12).Value = 4198.5714
$ws.Cells.Item(134, 13).Value = -4198.459800000001
$ws.Cells.Item(134, 14).Value = -9268.571400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 42111.11
$ws.Cells.Item(37, 10).Value = 42111.11
$ws.Cells.Item(37, 12).Value = 126333.33
$ws.Cells.Item(37, 14).Value = -126557.33

$ws.Cells.Item(92, 8).Value = 750.75
$ws.Cells.Item(92, 9).Value = 500
$ws.Cells.Item(92, 10).Value = 834.3333
$ws.Cells.Item(92, 11).Value = 1500
$ws.Cells.Item(92, 12).Value = 2502.9999
$ws.Cells.Item(92, 13).Value = -252
$ws.Cells.Item(92, 14).Value = -4998.9999

$ws.Cells.Item(117, 8).Value = 18527650
$ws.Cells.Item(117, 9).Value = 17038.166
$ws.Cells.Item(117, 10).Value = 27782954
$ws.Cells.Item(117, 11).Value = 51114.49800000001
$ws.Cells.Item(117, 12).Value = 83348862
$ws.Cells.Item(117, 13).Value = -47672.49800000001
$ws.Cells.Item(117, 14).Value = -83355746

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 33334360
$ws.Cells.Item(113, 9).Value = 58824344
$ws.Cells.Item(113, 10).Value = 1306.9231
$ws.Cells.Item(113, 11).Value = 58824344
$ws.Cells.Item(113, 12).Value = 1306.9231
$ws.Cells.Item(113, 13).Value = -58822174
$ws.Cells.Item(113, 14).Value = -5646.9231

$ws.Cells.Item(122, 8).Value = 24764532
$ws.Cells.Item(122, 9).Value = 39438744
$ws.Cells.Item(122, 10).Value = 1802.4375
$ws.Cells.Item(122, 11).Value = 118316232
$ws.Cells.Item(122, 12).Value = 5407.3125
$ws.Cells.Item(122, 13).Value = -118313782
$ws.Cells.Item(122, 14).Value = -10307.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 11877868
$ws.Cells.Item(132, 9).Value = 15719904
$ws.Cells.Item(132, 10).Value = 2481.6365
$ws.Cells.Item(132, 11).Value = 47159712
$ws.Cells.Item(132, 12).Value = 7444.9095
$ws.Cells.Item(132, 13).Value = -47157182
$ws.Cells.Item(132, 14).Value = -12504.9095

$ws.Cells.Item(136, 8).Value = 6306.5615
$ws.Cells.Item(136, 9).Value = 4182.9536
$ws.Cells.Item(136, 10).Value = 12829.071
$ws.Cells.Item(136, 11).Value = 12548.8608
$ws.Cells.Item(136, 12).Value = 38487.213
$ws.Cells.Item(136, 13).Value = -9998.860799999999
$ws.Cells.Item(136, 14).Value = -43587.213

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2961.524
$ws.Cells.Item(122, 9).Value = 2874.5
$ws.Cells.Item(122, 10).Value = 3240
$ws.Cells.Item(122, 11).Value = 8623.5
$ws.Cells.Item(122, 12).Value = 9720
$ws.Cells.Item(122, 13).Value = -6173.5
$ws.Cells.Item(122, 14).Value = -14620

$ws.Cells.Item(132, 8).Value = 16803.604
$ws.Cells.Item(132, 9).Value = 19555.924
$ws.Cells.Item(132, 10).Value = 2216.3
$ws.Cells.Item(132, 11).Value = 58667.772
$ws.Cells.Item(132, 12).Value = 6648.900000000001
$ws.Cells.Item(132, 13).Value = -56137.772
$ws.Cells.Item(132, 14).Value = -11708.9

$ws.Cells.Item(136, 8).Value = 8623339
$ws.Cells.Item(136, 9).Value = 2847.4866
$ws.Cells.Item(136, 11).Value = 8542.459800000001
$ws.Cells.Item(136, 13).Value = -5992.459800000001
